$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.852.43'
$ws.Range('E2').Value = '  +1.07%  '
$ws.Range('D3').Value = '3.543.74'
$ws.Range('E3').Value = '  +4.36%  '
$ws.Range('E4').Value = '  -0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '598.40'
$ws.Range('E5').Value = '  +3.61%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '136.41'
$ws.Range('E6').Value = '  +4.17%  '
$ws.Range('D7').Value = '3.541.34'
$ws.Range('E7').Value = '  +4.46%  '
$ws.Range('E8').Value = '  -0.08%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.496'
$ws.Range('E9').Value = '  +3.71%  '
$ws.Range('E10').Value = '  +3.05%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.91'
$ws.Range('E11').Value = '  -0.76%  '
$ws.Range('E12').Value = '  +4.41%  '
$ws.Range('D13').Value = '4.141.50'
$ws.Range('E13').Value = '  +4.11%  '
$ws.Range('E14').Value = '  +3.66%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '27.14'
$ws.Range('E15').Value = '  +5.46%  '
$ws.Range('D16').Value = '3.550.53'
$ws.Range('E16').Value = '  +5.89%  '
$ws.Range('E17').Value = '  +1.08%  '
$ws.Range('D18').Value = '64.762.90'
$ws.Range('E18').Value = '  +0.56%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '10.03'
$ws.Range('E19').Value = '  +5.82%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '14.39'
$ws.Range('E20').Value = '  +6.76%  '
$ws.Range('E21').Value = '  +3.67%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '388.38'
$ws.Range('E22').Value = '  +2.94%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.576'
$ws.Range('E23').Value = '  +7.01%  '
$ws.Range('D24').Value = '3.685.79'
$ws.Range('E24').Value = '  +4.35%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '74.24'
$ws.Range('E25').Value = '  +4.29%  '
$ws.Range('E26').Value = '  +0.01%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.0000116'
$ws.Range('E27').Value = '  +13.82%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.65'
$ws.Range('E28').Value = '  +8.77%  '
$ws.Range('E29').Value = '  +0.04%  '
$ws.Range('E30').Value = '  +5.62%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '8.34'
$ws.Range('E31').Value = '  +5.34%  '
$ws.Range('D32').Value = '3.549.76'
$ws.Range('E32').Value = '  +4.01%  '
$ws.Range('E33').Value = '  +24.50%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '23.95'
$ws.Range('E34').Value = '  +5.01%  '
$ws.Range('E35').Value = '  -0.01%  '
$ws.Range('E36').Value = '  +3.20%  '
$ws.Range('E37').Value = '  +5.45%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '169.74'
$ws.Range('E38').Value = '  -0.61%  '
$ws.Range('E39').Value = '  +7.24%  '
$ws.Range('E40').Value = '  +9.94%  '
$ws.Range('E41').Value = '  +7.13%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.825'
$ws.Range('E42').Value = '  +4.32%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '26.71'
$ws.Range('E43').Value = '  +21.58%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '42.61'
$ws.Range('E44').Value = '  +2.71%  '
$ws.Range('E45').Value = '  -0.25%  '
$ws.Range('E46').Value = '  +5.09%  '
$ws.Range('E47').Value = '  +10.36%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.64'
$ws.Range('E48').Value = '  +3.93%  '
$ws.Range('E49').Value = '  +7.30%  '
$ws.Range('D50').Value = '2.448.10'
$ws.Range('E50').Value = '  +12.42%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.34'
$ws.Range('E51').Value = '  +16.10%  '
